# Add Immigration test case (and Dependents sheet) to OrangeHRM_Testdata.xlsx
# Mirrors the authored diff:
#   - PIM_Add_Employee: new employee identity (Sunil Jain / 9456 / Sunil_Jain / Sunil@20)
#   - Employee_Personal_Details: gender Female -> Male
#   - Employee_Contact_Details: workemail/otheremail swapped to Rohan@gmail.com / Jain@gmail.com
#   - Emergency_Contacts: G1 header loses its (redundant/no-op) fill flag
#   - New sheets: Dependents, Immigration

$wb = $excel.ActiveWorkbook

$wsLogin   = $wb.Worksheets.Item("Login")
$wsPim     = $wb.Worksheets.Item("PIM_Add_Employee")
$wsPers    = $wb.Worksheets.Item("Employee_Personal_Details")
$wsContact = $wb.Worksheets.Item("Employee_Contact_Details")
$wsEmerg   = $wb.Worksheets.Item("Emergency_Contacts")

# ---------------------------------------------------------------------------
# 1. PIM_Add_Employee - new employee record
# ---------------------------------------------------------------------------
$wsPim.Range("A2").Value = "Sunil"
$wsPim.Range("C2").Value = "Jain"
$wsPim.Range("D2").Value = 9456
$wsPim.Range("E2").Value = "Sunil_Jain"
$wsPim.Range("F2").Value = "Sunil@20"
$wsPim.Range("G2").Value = "Sunil@20"
$wsPim.Activate()
$wsPim.Range("G2").Select()

# ---------------------------------------------------------------------------
# 2. Employee_Personal_Details - gender change
# ---------------------------------------------------------------------------
$wsPers.Range("K2").Value = "Male"
$wsPers.Activate()
$wsPers.Range("K2").Select()

# ---------------------------------------------------------------------------
# 3. Employee_Contact_Details - updated emails
# ---------------------------------------------------------------------------
$wsContact.Range("J2").Value = "Rohan@gmail.com"
$wsContact.Range("K2").Value = "Jain@gmail.com"
$wsContact.Activate()
$wsContact.Range("K2").Select()

# ---------------------------------------------------------------------------
# 4. Emergency_Contacts - drop the stray fill flag on G1 (copy format from a
#    cell that already uses the no-fill variant of the same border/font combo)
# ---------------------------------------------------------------------------
$wsLogin.Range("D1").Copy()
$wsEmerg.Range("G1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$wsEmerg.Activate()
$wsEmerg.Range("F1:F2").Select()

# ---------------------------------------------------------------------------
# 5. New sheet: Dependents
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsDep = $wb.Worksheets.Add($null, $lastSheet)
$wsDep.Name = "Dependents"

# Pull the header/data formatting from existing cells that already carry the
# exact styles used on this sheet (s=11 header, s=1 plain-bold header, s=2 data)
$wsContact.Range("D1").Copy()
$wsDep.Range("A1:E1").PasteSpecial(-4122)
$wsLogin.Range("A1").Copy()
$wsDep.Range("F1").PasteSpecial(-4122)
$wsLogin.Range("A2").Copy()
$wsDep.Range("A2:F2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsDep.Range("A1").Value = "dependentsname"
$wsDep.Range("B1").Value = "relationship"
$wsDep.Range("C1").Value = "dobmonth"
$wsDep.Range("D1").Value = "dobyear"
$wsDep.Range("E1").Value = "dobdate"
$wsDep.Range("F1").Value = "comment"

$wsDep.Range("A2").Value = "Mayank"
$wsDep.Range("B2").Value = "Child"
$wsDep.Range("C2").Value = "September"
$wsDep.Range("D2").Value = 2002
$wsDep.Range("E2").Value = 18
$wsDep.Range("F2").Value = "He is the child of the Employee."

for ($i = 1; $i -le 6; $i++) {
    $wsDep.Columns.Item($i).AutoFit()
}

$wsDep.Activate()
$wsDep.Range("H2").Select()

# ---------------------------------------------------------------------------
# 6. New sheet: Immigration
# ---------------------------------------------------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsImm = $wb.Worksheets.Add($null, $lastSheet2)
$wsImm.Name = "Immigration"

$wsContact.Range("D1").Copy()
$wsImm.Range("A1:O1").PasteSpecial(-4122)
$wsLogin.Range("A2").Copy()
$wsImm.Range("A2:O2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsImm.Range("A1").Value = "document"
$wsImm.Range("B1").Value = "number"
$wsImm.Range("C1").Value = "issuedby"
$wsImm.Range("D1").Value = "issuedmonth"
$wsImm.Range("E1").Value = "issuedyear"
$wsImm.Range("F1").Value = "issueddate"
$wsImm.Range("G1").Value = "expirymonth"
$wsImm.Range("H1").Value = "expiryyear"
$wsImm.Range("I1").Value = "expirydate"
$wsImm.Range("J1").Value = "eligiblestatus"
$wsImm.Range("K1").Value = "reviewmonth"
$wsImm.Range("L1").Value = "reviewyear"
$wsImm.Range("M1").Value = "reviewdate"
$wsImm.Range("N1").Value = "immigrationcomments"
$wsImm.Range("O1").Value = "comment"

$wsImm.Range("A2").Value = "Passport"
$wsImm.Range("B2").Value = 9812345670
$wsImm.Range("C2").Value = "India"
$wsImm.Range("D2").Value = "June"
$wsImm.Range("E2").Value = 2020
$wsImm.Range("F2").Value = 1
$wsImm.Range("G2").Value = "July"
$wsImm.Range("H2").Value = 2035
$wsImm.Range("I2").Value = 7
$wsImm.Range("J2").Value = "Yes"
$wsImm.Range("K2").Value = "September"
$wsImm.Range("L2").Value = 2023
$wsImm.Range("M2").Value = 11
$wsImm.Range("N2").Value = "Employee is eligble for immigration"
$wsImm.Range("O2").Value = "This is the Proof of Immigration"

for ($i = 1; $i -le 15; $i++) {
    $wsImm.Columns.Item($i).AutoFit()
}

$wsImm.PageSetup.Orientation = 1  # xlPortrait

$wsImm.Activate()
$wsImm.Range("L2").Select()

# ---------------------------------------------------------------------------
# Restore original active sheet/selection
# ---------------------------------------------------------------------------
$wsPim.Activate()
$wsPim.Range("G2").Select()

Write-Output "edit complete"
